$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the touched range to Text format so numeric-looking strings
# (e.g. "1.00", "41.900.91") are preserved verbatim, matching the
# original inlineStr cell type instead of being coerced to numbers.
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '41.900.91'
$ws.Range('E2').Value = '  -1.23%  '
$ws.Range('D3').Value = '2.487.53'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '312.09'
$ws.Range('E5').Value = '  +0.39%  '
$ws.Range('D6').Value = '95.45'
$ws.Range('E6').Value = '  -3.11%  '
$ws.Range('D7').Value = '0.555'
$ws.Range('E7').Value = '  -1.98%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '0.511'
$ws.Range('E9').Value = '  -3.07%  '
$ws.Range('D10').Value = '34.21'
$ws.Range('E10').Value = '  -4.05%  '
$ws.Range('D11').Value = '0.0787'
$ws.Range('E11').Value = '  -1.68%  '
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '7.05'
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('D14').Value = '2.870.45'
$ws.Range('E14').Value = '  -1.75%  '
$ws.Range('D15').Value = '2.508.32'
$ws.Range('E15').Value = '  -5.03%  '
$ws.Range('D16').Value = '14.84'
$ws.Range('E16').Value = '  -5.44%  '
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  -4.22%  '
$ws.Range('D18').Value = '41.937.04'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '6.44'
$ws.Range('E19').Value = '  -5.05%  '
$ws.Range('D20').Value = '0.0₃0925'
$ws.Range('E20').Value = '  -2.50%  '
$ws.Range('D21').Value = '11.76'
$ws.Range('E21').Value = '  -3.52%  '
$ws.Range('D22').Value = '69.56'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').Value = '238.14'
$ws.Range('E23').Value = '  -2.30%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').Value = '2.81'
$ws.Range('E24').Value = '  -3.06%  '
$ws.Range('D25').Value = '1.95'
$ws.Range('E25').Value = '  -4.63%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('D27').Value = '24.88'
$ws.Range('E27').Value = '  -3.97%  '
$ws.Range('E28').Value = '  -4.05%  '
$ws.Range('D29').Value = '9.78'
$ws.Range('E29').Value = '  -3.32%  '
$ws.Range('D30').Value = '36.76'
$ws.Range('E30').Value = '  -6.19%  '
$ws.Range('D31').Value = '154.79'
$ws.Range('E31').Value = '  -1.81%  '
$ws.Range('D32').Value = '5.66'
$ws.Range('E32').Value = '  -0.64%  '
$ws.Range('E33').Value = '  +0.12%  '
$ws.Range('B34').Value = 'ApeXProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  -8.33%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.0762'
$ws.Range('E35').Value = '  -4.14%  '
$ws.Range('D36').Value = '3.04'
$ws.Range('E36').Value = '  -3.54%  '
$ws.Range('D37').Value = '17.38'
$ws.Range('E37').Value = '  -3.80%  '
$ws.Range('D38').Value = '1.90'
$ws.Range('E38').Value = '  -5.81%  '
$ws.Range('D39').Value = '0.107'
$ws.Range('E39').Value = '  -3.58%  '
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('D41').Value = '4.04'
$ws.Range('E41').Value = '  -5.17%  '
$ws.Range('D42').Value = '21.44'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D44').Value = '2.003.96'
$ws.Range('E44').Value = '  +2.14%  '
$ws.Range('D45').Value = '0.0288'
$ws.Range('E45').Value = '  -2.89%  '
$ws.Range('D46').Value = '3.09'
$ws.Range('E46').Value = '  -5.95%  '
$ws.Range('D47').Value = '8.74'
$ws.Range('E47').Value = '  -2.04%  '
$ws.Range('D48').Value = '2.718.74'
$ws.Range('E48').Value = '  -2.06%  '
$ws.Range('D49').Value = '77.54'
$ws.Range('E49').Value = '  -4.29%  '
$ws.Range('D50').Value = '0.183'
$ws.Range('E50').Value = '  -4.25%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '98.36'
$ws.Range('E51').Value = '  -3.07%  '

# Restore the default cell style (no explicit style index), matching
# the original workbook where these data cells carried no "s" attribute.
$ws.Range('B2:E51').Style = 'Normal'
